$d = $word.ActiveDocument

# "Weather:  Sunny" -> "Authour: Monty"
$d.Content.Find.Execute("Weather:  Sunny", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Authour: Monty", 2)

# "Date: March 1" -> "Date: May 1"
$d.Content.Find.Execute("Date: March 1", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Date: May 1", 2)

# Append two new paragraphs at the end of the document: "end" and "Meep morp"
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Range.Text = "end"

$lastPara2 = $d.Paragraphs.Last
$lastPara2.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Range.Text = "Meep morp"
